$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column H with header "Hora de Reparacion" ---
$ws.Range("H1").Value = "Hora de Reparación"

# --- Touch row 2 (blank separator row) so it gets serialized as an empty row ---
# (toggling the outline level is a no-op but forces Excel to persist the row)
$ws.Rows("2").OutlineLevel = 1
$ws.Rows("2").OutlineLevel = 0

# --- Append new ticket rows 52-59 ---
# Row 52
$ws.Range("A52").Value = "'2024-05-15"
$ws.Range("B52").Value = "12:21:27"
$ws.Range("C52").Value = "-"
$ws.Range("D52").Value = "Cámara no detecta skeleton"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "12:21:37"

# Row 53
$ws.Range("A53").Value = "'2024-05-15"
$ws.Range("B53").Value = "12:21:56"
$ws.Range("C53").Value = "-"
$ws.Range("D53").Value = "Power atascado en prensa, cuesta sacar"
$ws.Range("E53").Value = "-"
$ws.Range("F53").Value = "-"
$ws.Range("G53").Value = "-"
$ws.Range("H53").Value = "12:22:17"

# Row 54
$ws.Range("A54").Value = "'2024-05-15"
$ws.Range("B54").Value = "12:22:37"
$ws.Range("C54").Value = "-"
$ws.Range("D54").Value = "AOI (fallo etiqueta)"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "12:22:42"

# Row 55
$ws.Range("A55").Value = "'2024-05-15"
$ws.Range("B55").Value = "12:22:57"
$ws.Range("C55").Value = "-"
$ws.Range("D55").Value = "Etiquetadora"
$ws.Range("E55").Value = "-"
$ws.Range("F55").Value = "-"
$ws.Range("G55").Value = "-"
$ws.Range("H55").Value = "12:23:36"

# Row 56
$ws.Range("A56").Value = "'2024-05-15"
$ws.Range("B56").Value = "12:23:54"
$ws.Range("C56").Value = "-"
$ws.Range("D56").Value = "Etiquetadora"
$ws.Range("E56").Value = "-"
$ws.Range("F56").Value = "-"
$ws.Range("G56").Value = "-"
$ws.Range("H56").Value = "12:23:58"

# Row 57
$ws.Range("A57").Value = "'2024-05-15"
$ws.Range("B57").Value = "12:24:03"
$ws.Range("C57").Value = "-"
$ws.Range("D57").Value = "Fallo etiqueta"
$ws.Range("E57").Value = "-"
$ws.Range("F57").Value = "-"
$ws.Range("G57").Value = "-"
$ws.Range("H57").Value = "12:24:05"

# Row 58
$ws.Range("A58").Value = "'2024-05-15"
$ws.Range("B58").Value = "12:24:17"
$ws.Range("C58").Value = "-"
$ws.Range("D58").Value = "Cámara no detecta busbar"
$ws.Range("E58").Value = "-"
$ws.Range("F58").Value = "-"
$ws.Range("G58").Value = "-"
$ws.Range("H58").Value = "12:24:28"

# Row 59
$ws.Range("A59").Value = "'2024-05-15"
$ws.Range("B59").Value = "12:25:51"
$ws.Range("C59").Value = "-"
$ws.Range("D59").Value = "Cámara no detecta Top cover"
$ws.Range("E59").Value = "-"
$ws.Range("F59").Value = "-"
$ws.Range("G59").Value = "-"
$ws.Range("H59").Value = "12:26:01"

# --- Normalize style on column A so the forced-text apostrophe does not leave a
#     quote-prefix cell format behind ---
$ws.Range("A52:A59").Style = "Normal"
